$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 766310
$ws.Range("E2").Value = 1429187466
$ws.Range("C10").Value = 345533
$ws.Range("E10").Value = 1817711330
$ws.Range("C13").Value = 187834
$ws.Range("D13").Value = 33262
$ws.Range("E13").Value = 1165224910
$ws.Range("C54").Value = 75191
$ws.Range("E54").Value = 361048094
$ws.Range("C78").Value = 178440
$ws.Range("E78").Value = 892473543
$ws.Range("C81").Value = 88350
$ws.Range("E81").Value = 499614738
$ws.Range("C88").Value = 71264
$ws.Range("E88").Value = 110294557
$ws.Range("C91").Value = 18848
$ws.Range("E91").Value = 75117834
$ws.Range("C93").Value = 16923
$ws.Range("E93").Value = 50452265
$ws.Range("C121").Value = 1306126
$ws.Range("D121").Value = 220385
$ws.Range("E121").Value = 2274557857
$ws.Range("C122").Value = 364
$ws.Range("E122").Value = 1156270
$ws.Range("C129").Value = 633338
$ws.Range("E129").Value = 3426710411
$ws.Range("C132").Value = 585622
$ws.Range("E132").Value = 3461185017
$ws.Range("C136").Value = 26675
$ws.Range("E136").Value = 143540532
$ws.Range("C139").Value = 76639
$ws.Range("E139").Value = 114131597
$ws.Range("C144").Value = 25068
$ws.Range("E144").Value = 92362703
$ws.Range("C151").Value = 39921
$ws.Range("E151").Value = 60359605
$ws.Range("C154").Value = 18438
$ws.Range("E154").Value = 72659539
$ws.Range("C156").Value = 12397
$ws.Range("E156").Value = 40042868
$ws.Range("C158").Value = 717
$ws.Range("E158").Value = 1762994
$ws.Range("C159").Value = 43847
$ws.Range("E159").Value = 101312981
$ws.Range("C178").Value = 515876
$ws.Range("E178").Value = 891189200
$ws.Range("C207").Value = 154660
$ws.Range("E207").Value = 753642756
